$wb = $excel.ActiveWorkbook

# OLE/VBA color ints for the workbook's existing "HyperLink" font color (RGB FF6495ED,
# stored BGR-swapped as 0x00BBGGRR for the Font.Color COM property).
$hyperlinkColor = 15570276  # 0xED9564  == RGB(0x64,0x95,0xED)

$statusHandedBack = "Handed back: in sync with en-US"

# Per-locale sheet configuration: handback timestamp + hyperlink target URLs
# (targets mirror the existing handoff-file / md-file hyperlinks already on
# each row, since the new "Latest Target File" / "Latest Handback File"
# columns reference the very same source .md and .xlf files).
$sheets = @(
    @{
        Name = "zh-cn"
        HandbackTime = "2016-02-22 09:21:38"
        Rows = @(
            @{
                Row = 2
                MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/f3d61203827af6255efaec5cb3502582b4fa82fd/e2e/04450948-5d02-4217-974e-d0ffa3ee09ff.md"
                MdDisplay = "04450948-5d02-4217-974e-d0ffa3ee09ff.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb7b8ad310fcaec6be65f9c6d53abe7f342a2b76/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/04450948-5d02-4217-974e-d0ffa3ee09ff.15a1edfc4b33e9a11595900e9d2ee44eec10b639.zh-cn.xlf"
                XlfDisplay = "04450948-5d02-4217-974e-d0ffa3ee09ff.15a1edfc4b33e9a11595900e9d2ee44eec10b639.zh-cn.xlf"
            },
            @{
                Row = 3
                MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/f3d61203827af6255efaec5cb3502582b4fa82fd/e2e/9f4ad892-dc37-4aec-b18d-4595c510be76.md"
                MdDisplay = "9f4ad892-dc37-4aec-b18d-4595c510be76.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb7b8ad310fcaec6be65f9c6d53abe7f342a2b76/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/9f4ad892-dc37-4aec-b18d-4595c510be76.1604df53df5c7b4bd28476df8b64ab4beae88868.zh-cn.xlf"
                XlfDisplay = "9f4ad892-dc37-4aec-b18d-4595c510be76.1604df53df5c7b4bd28476df8b64ab4beae88868.zh-cn.xlf"
            }
        )
    },
    @{
        Name = "de-de"
        HandbackTime = "2016-02-22 09:22:00"
        Rows = @(
            @{
                Row = 2
                MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/f3d61203827af6255efaec5cb3502582b4fa82fd/e2e/04450948-5d02-4217-974e-d0ffa3ee09ff.md"
                MdDisplay = "04450948-5d02-4217-974e-d0ffa3ee09ff.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9388be0a6abc8e023f8e9a4719d2f1c507264a94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/04450948-5d02-4217-974e-d0ffa3ee09ff.15a1edfc4b33e9a11595900e9d2ee44eec10b639.de-de.xlf"
                XlfDisplay = "04450948-5d02-4217-974e-d0ffa3ee09ff.15a1edfc4b33e9a11595900e9d2ee44eec10b639.de-de.xlf"
            },
            @{
                Row = 3
                MdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/f3d61203827af6255efaec5cb3502582b4fa82fd/e2e/9f4ad892-dc37-4aec-b18d-4595c510be76.md"
                MdDisplay = "9f4ad892-dc37-4aec-b18d-4595c510be76.md"
                XlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9388be0a6abc8e023f8e9a4719d2f1c507264a94/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/9f4ad892-dc37-4aec-b18d-4595c510be76.1604df53df5c7b4bd28476df8b64ab4beae88868.de-de.xlf"
                XlfDisplay = "9f4ad892-dc37-4aec-b18d-4595c510be76.1604df53df5c7b4bd28476df8b64ab4beae88868.de-de.xlf"
            }
        )
    }
)

foreach ($sheetCfg in $sheets) {
    $ws = $wb.Worksheets.Item($sheetCfg.Name)

    foreach ($rowCfg in $sheetCfg.Rows) {
        $r = $rowCfg.Row

        # Status column (B): handoff report generated for handback.
        $ws.Range("B$r").Value = $statusHandedBack

        # Latest Target File (E): hyperlink to the source .md file, same
        # target as the existing column-A hyperlink on this row.
        $eCell = $ws.Range("E$r")
        $ws.Hyperlinks.Add($eCell, $rowCfg.MdTarget, "", "", $rowCfg.MdDisplay)
        $eCell.Font.Underline = $true
        $eCell.Font.Color = $hyperlinkColor

        # Latest Handback File (F): hyperlink to the handed-back .xlf file,
        # same target as the existing column-C hyperlink on this row.
        $fCell = $ws.Range("F$r")
        $ws.Hyperlinks.Add($fCell, $rowCfg.XlfTarget, "", "", $rowCfg.XlfDisplay)
        $fCell.Font.Underline = $true
        $fCell.Font.Color = $hyperlinkColor

        # Latest Handback DateTime (G): stamp the handback time.
        $ws.Range("G$r").Value = $sheetCfg.HandbackTime

        # Handoff Reason (H): now included in the handback.
        $ws.Range("H$r").Value = "Include"
    }
}

Write-Output "Generated handback report for $($sheets.Count) locale sheets"
